$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 7.193183327378438
$ws.Cells.Item(2, 5).Value = 9.591339540850875
$ws.Cells.Item(3, 3).Value = -14.96173956806345
$ws.Cells.Item(3, 5).Value = -4.932343798304595
$ws.Cells.Item(4, 3).Value = 2.682935444832424
$ws.Cells.Item(4, 5).Value = -2.225127715916653
$ws.Cells.Item(5, 3).Value = 9.399485634179229
$ws.Cells.Item(5, 5).Value = 1.811802132286955
$ws.Cells.Item(6, 3).Value = 5.169490031659674
$ws.Cells.Item(6, 5).Value = 9.213376886330305
$ws.Cells.Item(7, 3).Value = -0.3722371047999662
$ws.Cells.Item(7, 5).Value = 2.684220738731935
$ws.Cells.Item(8, 3).Value = 4.098801479368341
$ws.Cells.Item(8, 5).Value = 2.548306621254004
$ws.Cells.Item(9, 3).Value = 3.75051862559701
$ws.Cells.Item(9, 5).Value = 2.714258593289975
$ws.Cells.Item(10, 3).Value = 2.352205130086071
$ws.Cells.Item(10, 5).Value = 3.873414041014778
$ws.Cells.Item(11, 3).Value = 4.083548352538369
$ws.Cells.Item(11, 5).Value = 3.586256146074462
$ws.Cells.Item(12, 3).Value = 4.861590900330692
$ws.Cells.Item(12, 5).Value = 3.297472770389764
$ws.Cells.Item(13, 3).Value = 1.787861866846807
$ws.Cells.Item(13, 5).Value = 4.088367525047842
$ws.Cells.Item(14, 3).Value = -2.21482332957591
$ws.Cells.Item(14, 5).Value = -0.6322362079330346
$ws.Cells.Item(15, 3).Value = 6.09521976277807
$ws.Cells.Item(15, 5).Value = 1.839905110456375
$ws.Cells.Item(16, 3).Value = 3.616930127707629
$ws.Cells.Item(16, 5).Value = 1.391416039405691
$ws.Cells.Item(17, 3).Value = 0.7171092762090492
$ws.Cells.Item(17, 5).Value = 2.755142438739822
$ws.Cells.Item(18, 3).Value = -0.1521036778360019
$ws.Cells.Item(18, 5).Value = 1.645968204809645
$ws.Cells.Item(19, 3).Value = -2.051528019634985
$ws.Cells.Item(19, 5).Value = -0.3224191428759626
